# Auto-generated: update leve-profit market-data columns (H-N) per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(21, 8).Value = 3682.2727
$ws.Cells.Item(21, 9).Value = 4003.4
$ws.Cells.Item(21, 10).Value = 3414.6667
$ws.Cells.Item(21, 11).Value = 4003.4
$ws.Cells.Item(21, 12).Value = 3414.6667
$ws.Cells.Item(21, 13).Value = -3535.4
$ws.Cells.Item(21, 14).Value = -4350.6667

$ws.Cells.Item(23, 8).Value = 3682.2727
$ws.Cells.Item(23, 9).Value = 4003.4
$ws.Cells.Item(23, 10).Value = 3414.6667
$ws.Cells.Item(23, 11).Value = 4003.4
$ws.Cells.Item(23, 12).Value = 3414.6667
$ws.Cells.Item(23, 13).Value = -3769.4
$ws.Cells.Item(23, 14).Value = -3882.6667

$ws.Cells.Item(29, 8).Value = 1779.7222
$ws.Cells.Item(29, 10).Value = 2856
$ws.Cells.Item(29, 12).Value = 8568
$ws.Cells.Item(29, 14).Value = -9130

$ws.Cells.Item(38, 8).Value = 240.46153
$ws.Cells.Item(38, 9).Value = 177.16667
$ws.Cells.Item(38, 10).Value = 1000
$ws.Cells.Item(38, 11).Value = 531.50001
$ws.Cells.Item(38, 12).Value = 3000
$ws.Cells.Item(38, 13).Value = -159.50001
$ws.Cells.Item(38, 14).Value = -3744

$ws.Cells.Item(86, 8).Value = 59540.5
$ws.Cells.Item(86, 9).Value = 4857.25
$ws.Cells.Item(86, 10).Value = 103287.1
$ws.Cells.Item(86, 11).Value = 4857.25
$ws.Cells.Item(86, 12).Value = 103287.1
$ws.Cells.Item(86, 13).Value = -3734.25
$ws.Cells.Item(86, 14).Value = -105533.1

$ws.Cells.Item(89, 8).Value = 59540.5
$ws.Cells.Item(89, 9).Value = 4857.25
$ws.Cells.Item(89, 10).Value = 103287.1
$ws.Cells.Item(89, 11).Value = 24286.25
$ws.Cells.Item(89, 12).Value = 516435.5
$ws.Cells.Item(89, 13).Value = -18670.25
$ws.Cells.Item(89, 14).Value = -527667.5

$ws.Cells.Item(125, 8).Value = 8751.556
$ws.Cells.Item(125, 9).Value = 8423.25
$ws.Cells.Item(125, 10).Value = 9014.200000000001
$ws.Cells.Item(125, 11).Value = 75809.25
$ws.Cells.Item(125, 12).Value = 81127.8
$ws.Cells.Item(125, 13).Value = -73349.25
$ws.Cells.Item(125, 14).Value = -86047.8

$ws.Cells.Item(138, 8).Value = 3308.2903
$ws.Cells.Item(138, 9).Value = 3876.1538
$ws.Cells.Item(138, 10).Value = 2898.1667
$ws.Cells.Item(138, 11).Value = 11628.4614
$ws.Cells.Item(138, 12).Value = 8694.500100000001
$ws.Cells.Item(138, 13).Value = -6488.4614
$ws.Cells.Item(138, 14).Value = -18974.5001


$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1042.75
$ws.Cells.Item(2, 9).Value = 939.7368
$ws.Cells.Item(2, 10).Value = 3000
$ws.Cells.Item(2, 11).Value = 939.7368
$ws.Cells.Item(2, 12).Value = 3000
$ws.Cells.Item(2, 13).Value = -826.7368
$ws.Cells.Item(2, 14).Value = -3226

$ws.Cells.Item(97, 8).Value = 597.3125
$ws.Cells.Item(97, 9).Value = 503.8
$ws.Cells.Item(97, 10).Value = 2000
$ws.Cells.Item(97, 11).Value = 503.8
$ws.Cells.Item(97, 12).Value = 2000
$ws.Cells.Item(97, 13).Value = -7.800000000000011
$ws.Cells.Item(97, 14).Value = -2992

$ws.Cells.Item(116, 8).Value = 1042.75
$ws.Cells.Item(116, 9).Value = 939.7368
$ws.Cells.Item(116, 10).Value = 3000
$ws.Cells.Item(116, 11).Value = 939.7368
$ws.Cells.Item(116, 12).Value = 3000
$ws.Cells.Item(116, 13).Value = 1354.2632
$ws.Cells.Item(116, 14).Value = -7588


$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1042.75
$ws.Cells.Item(3, 9).Value = 939.7368
$ws.Cells.Item(3, 10).Value = 3000
$ws.Cells.Item(3, 11).Value = 939.7368
$ws.Cells.Item(3, 12).Value = 3000
$ws.Cells.Item(3, 13).Value = -825.7368
$ws.Cells.Item(3, 14).Value = -3228

$ws.Cells.Item(7, 8).Value = 5001103.5
$ws.Cells.Item(7, 10).Value = 1493.6666
$ws.Cells.Item(7, 12).Value = 1493.6666
$ws.Cells.Item(7, 14).Value = -1719.6666

$ws.Cells.Item(94, 8).Value = 132.83333
$ws.Cells.Item(94, 9).Value = 132.83333
$ws.Cells.Item(94, 11).Value = 132.83333
$ws.Cells.Item(94, 13).Value = 318.16667


$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(13, 8).Value = 143.5
$ws.Cells.Item(13, 10).Value = 143.5
$ws.Cells.Item(13, 12).Value = 143.5
$ws.Cells.Item(13, 14).Value = -421.5

$ws.Cells.Item(31, 8).Value = 6919
$ws.Cells.Item(31, 9).Value = 4687.3335
$ws.Cells.Item(31, 10).Value = 9597
$ws.Cells.Item(31, 11).Value = 4687.3335
$ws.Cells.Item(31, 12).Value = 9597
$ws.Cells.Item(31, 13).Value = -4392.3335
$ws.Cells.Item(31, 14).Value = -10187

$ws.Cells.Item(34, 8).Value = 6919
$ws.Cells.Item(34, 9).Value = 4687.3335
$ws.Cells.Item(34, 10).Value = 9597
$ws.Cells.Item(34, 11).Value = 4687.3335
$ws.Cells.Item(34, 12).Value = 9597
$ws.Cells.Item(34, 13).Value = -4485.3335
$ws.Cells.Item(34, 14).Value = -10001


$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 111421.664
$ws.Cells.Item(107, 9).Value = 317.66666
$ws.Cells.Item(107, 10).Value = 166973.67
$ws.Cells.Item(107, 11).Value = 952.9999799999999
$ws.Cells.Item(107, 12).Value = 500921.01
$ws.Cells.Item(107, 13).Value = 967.0000200000001
$ws.Cells.Item(107, 14).Value = -504761.01


$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(92, 8).Value = 11080.167
$ws.Cells.Item(92, 10).Value = 11080.167
$ws.Cells.Item(92, 12).Value = 11080.167
$ws.Cells.Item(92, 14).Value = -14824.167

$ws.Cells.Item(97, 8).Value = 727.2105
$ws.Cells.Item(97, 9).Value = 632.93335
$ws.Cells.Item(97, 11).Value = 632.93335
$ws.Cells.Item(97, 13).Value = -136.93335


$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 11398.5
$ws.Cells.Item(16, 9).Value = 9899
$ws.Cells.Item(16, 10).Value = 12898
$ws.Cells.Item(16, 11).Value = 9899
$ws.Cells.Item(16, 12).Value = 12898
$ws.Cells.Item(16, 13).Value = -9729
$ws.Cells.Item(16, 14).Value = -13238

$ws.Cells.Item(40, 8).Value = 1937.75
$ws.Cells.Item(40, 9).Value = 1937.75
$ws.Cells.Item(40, 11).Value = 1937.75
$ws.Cells.Item(40, 13).Value = -1801.75

$ws.Cells.Item(82, 8).Value = 101089
$ws.Cells.Item(82, 9).Value = 1377.2858
$ws.Cells.Item(82, 10).Value = 333749.66
$ws.Cells.Item(82, 11).Value = 1377.2858
$ws.Cells.Item(82, 12).Value = 333749.66
$ws.Cells.Item(82, 13).Value = -1016.2858
$ws.Cells.Item(82, 14).Value = -334471.66

$ws.Cells.Item(85, 8).Value = 101089
$ws.Cells.Item(85, 9).Value = 1377.2858
$ws.Cells.Item(85, 10).Value = 333749.66
$ws.Cells.Item(85, 11).Value = 1377.2858
$ws.Cells.Item(85, 12).Value = 333749.66
$ws.Cells.Item(85, 13).Value = -129.2858000000001
$ws.Cells.Item(85, 14).Value = -336245.66

$ws.Cells.Item(93, 8).Value = 7332.3335
$ws.Cells.Item(93, 9).Value = 7332.3335
$ws.Cells.Item(93, 11).Value = 7332.3335
$ws.Cells.Item(93, 13).Value = -6084.3335


$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 3999
$ws.Cells.Item(81, 9).Value = 3999
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 11).Value = 7998
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 13).Value = -6937
$ws.Cells.Item(81, 14).ClearContents()

$ws.Cells.Item(84, 8).Value = 3999
$ws.Cells.Item(84, 9).Value = 3999
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 39990
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = -34686
$ws.Cells.Item(84, 14).ClearContents()

$ws.Cells.Item(100, 8).Value = 2225
$ws.Cells.Item(100, 9).Value = 1848.5
$ws.Cells.Item(100, 11).Value = 3697
$ws.Cells.Item(100, 13).Value = -3156

$ws.Cells.Item(132, 8).Value = 2689.3333
$ws.Cells.Item(132, 9).Value = 2487.2
$ws.Cells.Item(132, 10).Value = 3700
$ws.Cells.Item(132, 11).Value = 7461.599999999999
$ws.Cells.Item(132, 12).Value = 11100
$ws.Cells.Item(132, 13).Value = -4931.599999999999
$ws.Cells.Item(132, 14).Value = -16160

